$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hyperlinks: keep only the first one (F2), repoint it to the new URL; delete the rest ---
$hyperlinks = $ws.Hyperlinks
while ($hyperlinks.Count() -gt 1) {
    $last = $null
    foreach ($lnk in $hyperlinks) { $last = $lnk }
    $last.Delete() | Out-Null
}
$first = $null
foreach ($lnk in $hyperlinks) { $first = $lnk }
$first.Address = "https://www.lancers.jp/work/detail/5393175"

# --- Remove old data rows 3:11, keeping header row 1 and data row 2 ---
$ws.Range("A3:H11").EntireRow.Delete() | Out-Null

# --- Refresh row 2 with the new/updated record ---
$ws.Range("A2").Value = "2025-09-14 06:24:07"
$ws.Range("B2").Value = "【業務委託】アプリ開発の継続的パートナ募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5393175"
$ws.Range("G2").Value = 93
$ws.Range("H2").Value = "◆開発 ◇アプリ"

# --- Column widths: B 46 -> 23, H 21 -> 12 (offset by 5/6 to account for the
#     character-width <-> stored-width padding conversion) ---
$ws.Columns.Item(2).ColumnWidth = 23 - 0.8333333333333334
$ws.Columns.Item(8).ColumnWidth = 12 - 0.8333333333333334
